# "Add files via upload" — add a new "第五周 周二" (Week 5, Tuesday) progress
# block to the bottom of the plan-tracking sheet, and fill in the
# "完成情况" (completion status) column for the previous block
# ("第五周 周一").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fill in the completion-status column (C) for the existing last
#    week block (rows 80-84, under the "2017.9.25 第五周 周一" header).
# ---------------------------------------------------------------------
$ws.Range("C80").Value = "完成"
$ws.Range("C81").Value = "完成"
$ws.Range("C82").Value = "完成"
$ws.Range("C83").Value = "完成"
$ws.Range("C84").Value = "完成"

# ---------------------------------------------------------------------
# 2) Build a brand-new week block in rows 88-97, re-using the same
#    layout/formatting as the most recent block (rows 78-87):
#      row 88      -> merged date header (A88:D88)
#      row 89      -> column headers (人员 / 计划任务 / 完成情况 / 备注)
#      rows 90-94  -> one row per team member
#      rows 95-97  -> merged "总结：" summary block (A95:D97)
# ---------------------------------------------------------------------

# Copy the formatting (styles, number formats, alignment, fonts) of the
# previous block down onto the new block first, so every cell already
# has the right look before we fill in values.
$ws.Range("A78:D87").Copy()
$ws.Range("A88").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Merge the header / summary rows exactly like the template block.
$ws.Range("A88:D88").Merge()
$ws.Range("A95:D97").Merge()

# -- Row 88: new week/date header --
$ws.Range("A88").Value = "日期：2017.9.26 第五周 周二"

# -- Row 89: repeat the column headers --
$ws.Range("A89").Value = "人员"
$ws.Range("B89").Value = "计划任务"
$ws.Range("C89").Value = "完成情况"
$ws.Range("D89").Value = "备注"

# -- Rows 90-94: one row per team member with their planned task --
$ws.Range("A90").Value = "伍圣和"
$ws.Range("B90").Value = "账户管理模块相关数据库表修改和完善，学习Hbuilder的使用教程，（包括新建项目和创建HTML页面）"

$ws.Range("A91").Value = "龙荣盛"
$ws.Range("B91").Value = "作品管理模块相关数据库表修改和完善，学习Hbuilder的使用教程，（包括新建项目和创建HTML页面）"

$ws.Range("A92").Value = "李志华"
$ws.Range("B92").Value = "资讯管理模块相关数据库表修改和完善，学习Hbuilder的使用教程，（包括新建项目和创建HTML页面）"

$ws.Range("A93").Value = "石婉霞"
$ws.Range("B93").Value = "布料管理模块相关数据库表修改和完善，学习Hbuilder的使用教程，（包括新建项目和创建HTML页面）"

$ws.Range("A94").Value = "陈俊彬"
$ws.Range("B94").Value = "服装管理模块相关数据库表修改和完善，学习Hbuilder的使用教程，（包括新建项目和创建HTML页面）"

# C90:D94 (completion status / remarks) are left blank for this new week,
# same as how the block started out before being filled in.

# -- Row 95-97: "总结：" summary placeholder, same text as other blocks --
$ws.Range("A95").Value = "总结："

# Row heights for the wrapped task-description rows, matching how Excel
# grows a row to fit newly typed, wrapped text.
$ws.Rows.Item(90).RowHeight = 69.75
$ws.Rows.Item(91).RowHeight = 71.25
$ws.Rows.Item(92).RowHeight = 66
$ws.Rows.Item(93).RowHeight = 63
$ws.Rows.Item(94).RowHeight = 66.75

# ---------------------------------------------------------------------
# 3) Scroll/selection bookkeeping so the view matches where the user
#    left off editing (best effort).
# ---------------------------------------------------------------------
$ws.Range("B90").Select()
